$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 389 (new week data replacing prior entry for Sutil De Gase)
$ws.Cells.Item(389, 4).Value = 45131
$ws.Cells.Item(389, 13).Value = 250
$ws.Cells.Item(389, 14).Value = 28000
$ws.Cells.Item(389, 15).Value = 29000
$ws.Cells.Item(389, 16).Value = 28400
$ws.Cells.Item(389, 19).Value = 1183

# Update existing row 390 (new week data replacing prior entry for Tahití / Colombia)
$ws.Cells.Item(390, 4).Value = 45131
$ws.Cells.Item(390, 13).Value = 400
$ws.Cells.Item(390, 14).Value = 26000
$ws.Cells.Item(390, 15).Value = 27000
$ws.Cells.Item(390, 16).Value = 26500
$ws.Cells.Item(390, 19).Value = 1104

# Insert 3 new rows at 391:393, pushing the rest of the table down
$ws.Range("A391:A393").EntireRow.Insert()

# Row 391
$ws.Cells.Item(391, 1).Value = 1
$ws.Cells.Item(391, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(391, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(391, 4).Value = 45131
$ws.Cells.Item(391, 5).Value = 15
$ws.Cells.Item(391, 6).Value = "Fruta"
$ws.Cells.Item(391, 7).Value = 100102
$ws.Cells.Item(391, 8).Value = "Cítricos"
$ws.Cells.Item(391, 9).Value = 100102003
$ws.Cells.Item(391, 10).Value = "Limón"
$ws.Cells.Item(391, 11).Value = "Tahití"
$ws.Cells.Item(391, 12).Value = "Primera"
$ws.Cells.Item(391, 13).Value = 500
$ws.Cells.Item(391, 14).Value = 28000
$ws.Cells.Item(391, 15).Value = 29000
$ws.Cells.Item(391, 16).Value = 28600
$ws.Cells.Item(391, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(391, 18).Value = "Perú"
$ws.Cells.Item(391, 19).Value = 1192
$ws.Cells.Item(391, 20).Value = 24

# Row 392
$ws.Cells.Item(392, 1).Value = 1
$ws.Cells.Item(392, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(392, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(392, 4).Value = 45124
$ws.Cells.Item(392, 5).Value = 15
$ws.Cells.Item(392, 6).Value = "Fruta"
$ws.Cells.Item(392, 7).Value = 100102
$ws.Cells.Item(392, 8).Value = "Cítricos"
$ws.Cells.Item(392, 9).Value = 100102003
$ws.Cells.Item(392, 10).Value = "Limón"
$ws.Cells.Item(392, 11).Value = "Sutil De Gase"
$ws.Cells.Item(392, 12).Value = "Primera"
$ws.Cells.Item(392, 13).Value = 300
$ws.Cells.Item(392, 14).Value = 40000
$ws.Cells.Item(392, 15).Value = 41000
$ws.Cells.Item(392, 16).Value = 40500
$ws.Cells.Item(392, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(392, 18).Value = "Perú"
$ws.Cells.Item(392, 19).Value = 1688
$ws.Cells.Item(392, 20).Value = 24

# Row 393
$ws.Cells.Item(393, 1).Value = 1
$ws.Cells.Item(393, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(393, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(393, 4).Value = 45124
$ws.Cells.Item(393, 5).Value = 15
$ws.Cells.Item(393, 6).Value = "Fruta"
$ws.Cells.Item(393, 7).Value = 100102
$ws.Cells.Item(393, 8).Value = "Cítricos"
$ws.Cells.Item(393, 9).Value = 100102003
$ws.Cells.Item(393, 10).Value = "Limón"
$ws.Cells.Item(393, 11).Value = "Tahití"
$ws.Cells.Item(393, 12).Value = "Primera"
$ws.Cells.Item(393, 13).Value = 350
$ws.Cells.Item(393, 14).Value = 28000
$ws.Cells.Item(393, 15).Value = 29000
$ws.Cells.Item(393, 16).Value = 28571
$ws.Cells.Item(393, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(393, 18).Value = "Colombia"
$ws.Cells.Item(393, 19).Value = 1190
$ws.Cells.Item(393, 20).Value = 24
